$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Justifikasi / bukti rencana usaha dan/atau kegiatan secara prinsip dapat dilakukan",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Justifikasi / bukti persetujuan awal rencana usaha dan/atau kegiatan",
    2
)
